$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("da_profit_df")

$headers = @("x1","x2","x3","x4","x5","x6","x7","x8","x9","x10","x11","x12","x13","x14","x15","x16","x17","x18","x19","x20","x21","x22","x23","x24")
$values = @(71865.48750000003,0.0,0.0,66867.11249999999,0.0,85071.02999999998,0.0,0.0,0.0,0.0,0.0,76175.91,71629.39499999999,69435.38250000002,0.0,0.0,89052.57,0.0,0.0,0.0,0.0,0.0,0.0,0.0)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headers[$i]
    $ws.Cells.Item(2, $col).Value = $values[$i]
}

# Remove old data below/around row 2 that is no longer part of the new layout
for ($r = 3; $r -le 25; $r++) {
    $ws.Cells.Item($r, 1).ClearContents()
}
